$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new product row (row 8: PROD7TEST / Product 7 test) mirroring the
# layout of the existing rows. Text-like values ("true"/"True" and the
# numeric-looking price strings) are entered with a leading apostrophe so
# Excel stores them as literal text (matching the shared-string/text cell
# type used throughout the rest of this import-template sample) instead of
# auto-converting them to booleans or numbers.
$ws.Range("A8").Value = "'PROD7TEST"
$ws.Range("B8").Value = "'Product 7 test"
$ws.Range("C8").Value = "'Categ 1"
$ws.Range("D8").Value = "'service"
$ws.Range("E8").Value = "'true"
$ws.Range("F8").Value = "'True"
$ws.Range("G8").Value = "'11.11"
$ws.Range("H8").Value = "'15.99"
# I8 (invoice_policy) is left blank for this row, same as the sheet's source.

# Move/refresh the active selection, matching the saved cursor position.
$ws.Range("F17").Select() | Out-Null
